# Atualização de bases das ligas, do dia: 14-04-2024 às 18:28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 114 / Row 115: the two fixtures (ids 7559468 / 7559469) had their
# data rows swapped (everything except the leading id-index column A and its
# style, which stay put). Row 114 now holds what used to be row 115's data
# and vice versa.

# Row 114 <- new values
$ws.Cells.Item(114, 2).Value  = 7559468
$ws.Cells.Item(114, 6).Value  = "Liverpool Montevideo"
$ws.Cells.Item(114, 7).Value  = "CA River Plate"
$ws.Cells.Item(114, 8).Value  = 2
$ws.Cells.Item(114, 9).Value  = 1
$ws.Cells.Item(114, 10).Value = "H"
$ws.Cells.Item(114, 11).Value = 1.7
$ws.Cells.Item(114, 12).Value = 3
$ws.Cells.Item(114, 13).Value = 5.75
$ws.Cells.Item(114, 14).Value = 1.833
$ws.Cells.Item(114, 15).Value = 3.2
$ws.Cells.Item(114, 16).Value = 4.5
$ws.Cells.Item(114, 17).Value = -0.5
$ws.Cells.Item(114, 18).Value = 1.925
$ws.Cells.Item(114, 19).Value = 1.925
$ws.Cells.Item(114, 20).Value = 2.25
$ws.Cells.Item(114, 21).Value = 2.025
$ws.Cells.Item(114, 22).Value = 1.825
$ws.Cells.Item(114, 23).Value = 0.833
$ws.Cells.Item(114, 24).Value = -1
$ws.Cells.Item(114, 25).Value = -1
$ws.Cells.Item(114, 26).Value = 0.925
$ws.Cells.Item(114, 27).Value = -1
$ws.Cells.Item(114, 28).Value = 1.025
$ws.Cells.Item(114, 29).Value = -1

# Row 115 <- new values
$ws.Cells.Item(115, 2).Value  = 7559469
$ws.Cells.Item(115, 6).Value  = "Montevideo Wanderers"
$ws.Cells.Item(115, 7).Value  = "Penarol"
$ws.Cells.Item(115, 8).Value  = 0
$ws.Cells.Item(115, 9).Value  = 0
$ws.Cells.Item(115, 10).Value = "D"
$ws.Cells.Item(115, 11).Value = 4.75
$ws.Cells.Item(115, 12).Value = 3.4
$ws.Cells.Item(115, 13).Value = 1.7
$ws.Cells.Item(115, 14).Value = 2.7
$ws.Cells.Item(115, 15).Value = 3.2
$ws.Cells.Item(115, 16).Value = 2.45
$ws.Cells.Item(115, 17).Value = 0
$ws.Cells.Item(115, 18).Value = 2.05
$ws.Cells.Item(115, 19).Value = 1.8
$ws.Cells.Item(115, 20).Value = 2.5
$ws.Cells.Item(115, 21).Value = 1.975
$ws.Cells.Item(115, 22).Value = 1.875
$ws.Cells.Item(115, 23).Value = -1
$ws.Cells.Item(115, 24).Value = 2.2
$ws.Cells.Item(115, 25).Value = -1
$ws.Cells.Item(115, 26).Value = 0
$ws.Cells.Item(115, 27).Value = -0
$ws.Cells.Item(115, 28).Value = -1
$ws.Cells.Item(115, 29).Value = 0.875

# --- Row 119 / Row 120: fixtures (ids 7013409 / 7013702) likewise swapped.

# Row 119 <- new values
$ws.Cells.Item(119, 2).Value  = 7013702
$ws.Cells.Item(119, 6).Value  = "Defensor Sporting"
$ws.Cells.Item(119, 7).Value  = "Danubio"
$ws.Cells.Item(119, 8).Value  = 0
$ws.Cells.Item(119, 9).Value  = 2
$ws.Cells.Item(119, 10).Value = "A"
$ws.Cells.Item(119, 11).Value = 1.8
$ws.Cells.Item(119, 12).Value = 3.6
$ws.Cells.Item(119, 13).Value = 4.2
$ws.Cells.Item(119, 14).Value = 1.8
$ws.Cells.Item(119, 15).Value = 3.6
$ws.Cells.Item(119, 16).Value = 4.2
$ws.Cells.Item(119, 17).Value = -0.75
$ws.Cells.Item(119, 18).Value = 2.05
$ws.Cells.Item(119, 19).Value = 1.8
$ws.Cells.Item(119, 20).Value = 2.25
$ws.Cells.Item(119, 21).Value = 1.85
$ws.Cells.Item(119, 22).Value = 2
$ws.Cells.Item(119, 23).Value = -1
$ws.Cells.Item(119, 24).Value = -1
$ws.Cells.Item(119, 25).Value = 3.2
$ws.Cells.Item(119, 26).Value = -1
$ws.Cells.Item(119, 27).Value = 0.8
$ws.Cells.Item(119, 28).Value = -0.5
$ws.Cells.Item(119, 29).Value = 0.5

# Row 120 <- new values
$ws.Cells.Item(120, 2).Value  = 7013409
$ws.Cells.Item(120, 6).Value  = "Nacional De Football"
$ws.Cells.Item(120, 7).Value  = "Torque"
$ws.Cells.Item(120, 8).Value  = 1
$ws.Cells.Item(120, 9).Value  = 1
$ws.Cells.Item(120, 10).Value = "D"
$ws.Cells.Item(120, 11).Value = 1.666
$ws.Cells.Item(120, 12).Value = 3.9
$ws.Cells.Item(120, 13).Value = 4.5
$ws.Cells.Item(120, 14).Value = 1.615
$ws.Cells.Item(120, 15).Value = 4
$ws.Cells.Item(120, 16).Value = 4.75
$ws.Cells.Item(120, 17).Value = -0.75
$ws.Cells.Item(120, 18).Value = 1.8
$ws.Cells.Item(120, 19).Value = 2.05
$ws.Cells.Item(120, 20).Value = 2.75
$ws.Cells.Item(120, 21).Value = 1.95
$ws.Cells.Item(120, 22).Value = 1.9
$ws.Cells.Item(120, 23).Value = -1
$ws.Cells.Item(120, 24).Value = 3
$ws.Cells.Item(120, 25).Value = -1
$ws.Cells.Item(120, 26).Value = -1
$ws.Cells.Item(120, 27).Value = 1.05
$ws.Cells.Item(120, 28).Value = -1
$ws.Cells.Item(120, 29).Value = 0.8999999999999999

# --- Rows 179-182 are removed entirely: matches 8051004 and 8051187 drop
# out of the sheet, while matches 8050911 and 8050912 (formerly rows 181
# and 182) move up into rows 177/178 -- replacing the old 8051185/8051186
# rows there -- carrying refreshed (closing) odds values.
$ws.Rows("179:182").Delete()

# Row 177 <- refreshed values for match id 8050911
$ws.Cells.Item(177, 2).Value  = 8050911
$ws.Cells.Item(177, 5).Value  = 45396.75
$ws.Cells.Item(177, 6).Value  = "Penarol"
$ws.Cells.Item(177, 7).Value  = "Danubio"
$ws.Cells.Item(177, 11).Value = 1.666
$ws.Cells.Item(177, 12).Value = 3.5
$ws.Cells.Item(177, 13).Value = 5.5
$ws.Cells.Item(177, 14).Value = 1.6
$ws.Cells.Item(177, 15).Value = 3.6
$ws.Cells.Item(177, 16).Value = 6.5
$ws.Cells.Item(177, 17).Value = -1
$ws.Cells.Item(177, 18).Value = 2.1
$ws.Cells.Item(177, 19).Value = 1.775
$ws.Cells.Item(177, 20).Value = 2.25
$ws.Cells.Item(177, 21).Value = 1.95
$ws.Cells.Item(177, 22).Value = 1.9

# Row 178 <- refreshed values for match id 8050912
$ws.Cells.Item(178, 2).Value  = 8050912
$ws.Cells.Item(178, 5).Value  = 45397.75
$ws.Cells.Item(178, 6).Value  = "Montevideo Wanderers"
$ws.Cells.Item(178, 7).Value  = "Liverpool Montevideo"
$ws.Cells.Item(178, 11).Value = 3.2
$ws.Cells.Item(178, 12).Value = 3.3
$ws.Cells.Item(178, 13).Value = 2.2
$ws.Cells.Item(178, 14).Value = 3.4
$ws.Cells.Item(178, 15).Value = 3.3
$ws.Cells.Item(178, 16).Value = 2.1
$ws.Cells.Item(178, 17).Value = 0.25
$ws.Cells.Item(178, 18).Value = 2.025
$ws.Cells.Item(178, 19).Value = 1.825
$ws.Cells.Item(178, 20).Value = 2.25
$ws.Cells.Item(178, 21).Value = 1.9
$ws.Cells.Item(178, 22).Value = 1.95
